$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("рис. fig.", $true, $false, $false, $false, $false, $true, 1, $false, "fig.", 2)
